$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the two target list-item paragraphs by their current text
# rather than hard-coded indices, so the script stays robust.
# ------------------------------------------------------------------
$queuePara = $null
$mapsPara  = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Improved queue management*") {
        $queuePara = $p
    } elseif ($t -like "*Google maps to show the position of your taxi*") {
        $mapsPara = $p
    }
}

if ($queuePara -eq $null) { throw "Could not find the '[Improved queue management?]' paragraph" }
if ($mapsPara -eq $null) { throw "Could not find the '[Google maps...]' paragraph" }

# ------------------------------------------------------------------
# Step 1: delete the pre-existing "_GoBack" bookmark (it currently
# trails the Google-maps paragraph); it will be re-created further
# up, inside the re-written queue-management paragraph. "_GoBack" is
# a hidden bookmark (leading underscore) so it will not show up via
# plain enumeration - look it up by name instead.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldGoBack = $d.Bookmarks.Item("_GoBack")
    $oldGoBack.Delete()
}

# ------------------------------------------------------------------
# Step 2: rewrite the "[Improved queue management?]" paragraph with
# the new sentence, split across several runs, with the _GoBack
# bookmark re-inserted after "with r" and proofErr gramStart/gramEnd
# markers wrapped around "high density".
# ------------------------------------------------------------------
$rPr = '<w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/></w:rPr>'

$queueXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r>$rPr<w:t>An improved taxi management system, with r</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r>$rPr<w:t xml:space="preserve">edistributes available taxis moving them from </w:t></w:r><w:proofErr w:type="gramStart"/><w:r>$rPr<w:t>high density</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r>$rPr<w:t xml:space="preserve"> areas to low density areas. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$queueRange = $queuePara.Range
$queueTarget = $d.Range($queueRange.Start, $queueRange.End)
$queueTarget.InsertXML($queueXml)

# ------------------------------------------------------------------
# Step 3: rewrite the "[Google maps to show the position of your
# taxi?]" paragraph text (the trailing bookmark was already removed
# in Step 1).
# ------------------------------------------------------------------
$mapsPara.Range.Text = "The possibility to show the GPS position of the customer’s taxi on an interactive map."
